# Update cryptos list values (price & 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.418.42"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.847.79"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'240.40"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'0.6333"
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.07560"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.2967"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "'24.60"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "'0.07717"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.847.95"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "'0.6859"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'0.00001008"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").Value = "'83.14"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "'6.179"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "29.437.27"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'230.22"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "'12.47"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'7.580"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'156.92"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "'0.1403"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").Value = "'8.381"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "'17.69"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'1.466"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'0.05723"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").Value = "'4.041"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "'1.852"
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").Value = "'1.158"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "'0.7170"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'2.594"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "1.252.17"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").Value = "'0.01819"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").Value = "'2.784"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.197"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9049"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").Value = "'0.9998"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "2.001.96"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "'101.81"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'66.41"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'9.182"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'7.062"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "'0.4036"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000117"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "'1.710"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "'0.1133"
$ws.Range("E51").Value = "  +1.23%  "
